{"js": "// Replace the 20x5 grid of arithmetic-practice cells with the new\n// problem set. Cell order in the source diff is row-major (top-to-bottom,\n// left-to-right), matching Word.Table.values, so we overwrite the whole\n// table in one positional assignment. This avoids ambiguity from the\n// duplicate / cross-referencing text values that appear across cells\n// (e.g. one cell's old text equals another cell's new text).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount, values\");\nawait context.sync();\n\nconst newValues = [\n  [\n    \"10+68=\",\n    \"2+96=\",\n    \"37+37=\",\n    \"15+18=\",\n    \"89-27=\"\n  ],\n  [\n    \"44-6=\",\n    \"74-59=\",\n    \"18+39=\",\n    \"99-72=\",\n    \"99-11=\"\n  ],\n  [\n    \"1+1=\",\n    \"91-15=\",\n    \"33+14=\",\n    \"81-69=\",\n    \"22+5=\"\n  ],\n  [\n    \"10+35=\",\n    \"53-27=\",\n    \"24+20=\",\n    \"66+12=\",\n    \"5+94=\"\n  ],\n  [\n    \"78+3=\",\n    \"17-15=\",\n    \"53-23=\",\n    \"28+40=\",\n    \"97-41=\"\n  ],\n  [\n    \"3+75=\",\n    \"79-32=\",\n    \"36-5=\",\n    \"57+19=\",\n    \"41-16=\"\n  ],\n  [\n    \"14+75=\",\n    \"21+77=\",\n    \"5+2=\",\n    \"48-29=\",\n    \"6+34=\"\n  ],\n  [\n    \"6+92=\",\n    \"72-47=\",\n    \"74+14=\",\n    \"17-15=\",\n    \"46-13=\"\n  ],\n  [\n    \"59-58=\",\n    \"72-68=\",\n    \"44+54=\",\n    \"73-43=\",\n    \"55+32=\"\n  ],\n  [\n    \"9+30=\",\n    \"58-36=\",\n    \"22+2=\",\n    \"34+18=\",\n    \"44-23=\"\n  ],\n  [\n    \"43-19=\",\n    \"66-37=\",\n    \"30+38=\",\n    \"7+57=\",\n    \"34+15=\"\n  ],\n  [\n    \"34+63=\",\n    \"25+23=\",\n    \"80+19=\",\n    \"21+65=\",\n    \"14+61=\"\n  ],\n  [\n    \"34-22=\",\n    \"88-40=\",\n    \"36-12=\",\n    \"20+21=\",\n    \"64+31=\"\n  ],\n  [\n    \"6+79=\",\n    \"19+48=\",\n    \"92-28=\",\n    \"16+57=\",\n    \"74+16=\"\n  ],\n  [\n    \"90-62=\",\n    \"40-34=\",\n    \"96-61=\",\n    \"42+25=\",\n    \"93-49=\"\n  ],\n  [\n    \"31-17=\",\n    \"75+16=\",\n    \"23+44=\",\n    \"65-5=\",\n    \"30+33=\"\n  ],\n  [\n    \"41-7=\",\n    \"93-11=\",\n    \"65+31=\",\n    \"5+77=\",\n    \"23+51=\"\n  ],\n  [\n    \"76+16=\",\n    \"17+51=\",\n    \"46-18=\",\n    \"86-69=\",\n    \"54+38=\"\n  ],\n  [\n    \"32+57=\",\n    \"83-82=\",\n    \"69-37=\",\n    \"33-9=\",\n    \"1+45=\"\n  ],\n  [\n    \"23-6=\",\n    \"59-37=\",\n    \"94-38=\",\n    \"26+23=\",\n    \"60+13=\"\n  ]\n];\n\nif (table.rowCount !== newValues.length) {\n  throw new Error(`Expected ${newValues.length} rows, found ${table.rowCount}.`);\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace the 20x5 grid of arithmetic-practice cells with the new problem\n# set. Cell order in the source diff is row-major (top-to-bottom,\n# left-to-right); Table.Cell(row, col) is 1-indexed, so we overwrite every\n# cell by position. Positional addressing (rather than text find/replace)\n# avoids ambiguity from duplicate / cross-referencing text values that\n# appear across cells (e.g. one cell's old text equals another cell's new\n# text, and one new value is shared by two different cells).\n$d = $word.ActiveDocument\n\nif ($d.Tables.Count -lt 1) {\n    throw \"Expected a table in the document body, found none.\"\n}\n\n$tbl = $d.Tables.Item(1)\n\n$newValues = @(\n    @('10+68=', '2+96=', '37+37=', '15+18=', '89-27='),\n    @('44-6=', '74-59=', '18+39=', '99-72=', '99-11='),\n    @('1+1=', '91-15=', '33+14=', '81-69=', '22+5='),\n    @('10+35=', '53-27=', '24+20=', '66+12=', '5+94='),\n    @('78+3=', '17-15=', '53-23=', '28+40=', '97-41='),\n    @('3+75=', '79-32=', '36-5=', '57+19=', '41-16='),\n    @('14+75=', '21+77=', '5+2=', '48-29=', '6+34='),\n    @('6+92=', '72-47=', '74+14=', '17-15=', '46-13='),\n    @('59-58=', '72-68=', '44+54=', '73-43=', '55+32='),\n    @('9+30=', '58-36=', '22+2=', '34+18=', '44-23='),\n    @('43-19=', '66-37=', '30+38=', '7+57=', '34+15='),\n    @('34+63=', '25+23=', '80+19=', '21+65=', '14+61='),\n    @('34-22=', '88-40=', '36-12=', '20+21=', '64+31='),\n    @('6+79=', '19+48=', '92-28=', '16+57=', '74+16='),\n    @('90-62=', '40-34=', '96-61=', '42+25=', '93-49='),\n    @('31-17=', '75+16=', '23+44=', '65-5=', '30+33='),\n    @('41-7=', '93-11=', '65+31=', '5+77=', '23+51='),\n    @('76+16=', '17+51=', '46-18=', '86-69=', '54+38='),\n    @('32+57=', '83-82=', '69-37=', '33-9=', '1+45='),\n    @('23-6=', '59-37=', '94-38=', '26+23=', '60+13='),\n)\n\nif ($tbl.Rows.Count -ne $newValues.Length) {\n    throw \"Expected $($newValues.Length) rows, found $($tbl.Rows.Count).\"\n}\n\nfor ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n    $rowValues = $newValues[$r - 1]\n    for ($c = 1; $c -le $tbl.Columns.Count; $c++) {\n        $tbl.Cell($r, $c).Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
